$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Remove the review comments that were left on B7/B8 (the "yellow" highlight
# discussion is resolved - Lizzie's email is applied below).
$ws1.Range("B7").Comment.Delete()
$ws1.Range("B8").Comment.Delete()

# Match the formatting used elsewhere in the sheet instead of the special
# "needs review / yellow" highlighting:
#  - B7 takes on the header-like black/size-14 style also used by B2
#  - B8 takes on the plain default row style used by the rest of column B
$ws1.Range("B2").Copy() | Out-Null
$ws1.Range("B7").PasteSpecial(-4122) | Out-Null

$ws1.Range("B9").Copy() | Out-Null
$ws1.Range("B8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Update the column mappings per Lizzie's email: use EMPLOYMENT_BEGIN_DATE /
# EMPLOYMENT_END_DATE instead of REQUESTED_BEGIN_DATE / REQUESTED_END_DATE.
$ws1.Range("B7").Value = "EMPLOYMENT_BEGIN_DATE"
$ws1.Range("B8").Value = "EMPLOYMENT_END_DATE"

# Move the active selection to B12, matching where the author left off editing.
$ws1.Range("B12").Select()
